$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tuesday (row 14): record the first in/out time pair (08:00 - 12:00)
$ws.Range("C14").Value = 0.333333333333333
$ws.Range("D14").Value = 0.5

# Free-text note describing the split shift worked that day
$ws.Range("L14").Value = "0800 – 1000, 1200 – 1400"

# D15 gets cleared out (content + the time formatting/validation that went
# with it) - copy the plain default style from a blank, unformatted cell
$ws.Range("D15").Validation.Delete()
$ws.Range("D15").ClearContents()
$ws.Range("B1").Copy()
$ws.Range("D15").PasteSpecial(-4122)

# Re-apply the print area (mirrors the repeated "Set Print Area" action
# that keeps stacking a fresh _xlnm.Print_Area_... name on save)
$wb.Names.Add("_xlnm.Print_Area_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0", "='Weekly Time Record'!`$A`$1:`$K`$27")

# Move the cursor to where it ended up after the edits
$ws.Range("D15").Select() | Out-Null
